# Updates the cryptos list worksheet with refreshed prices / volume figures,
# mirroring the daily GitHub Actions data-refresh commit.
#
# All of the Price/Volume/Coin/Link cells in this sheet are stored as text
# (inline strings) rather than numbers -- values such as "28.212.54" or
# "1.005" are thousand-grouped / percentage strings, not numerics. Plainly
# assigning $cell.Value would let Excel auto-coerce number-looking strings
# into floating point values (e.g. "1.005" -> 1.0049999999999999) and drop
# leading/trailing zeros. To avoid that we force the cell to text format
# before assigning the value, then restore the default "Normal" style so we
# don't leave stray formatting behind.

function Set-CellText($sheet, $ref, $val) {
    $c = $sheet.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 37/38 swapped position (VeChain now ranks above Aptos) ---
Set-CellText $ws "B37" "VeChain"
Set-CellText $ws "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D37" "0.02351"
Set-CellText $ws "E37" "  +4.51%  "

Set-CellText $ws "B38" "Aptos"
Set-CellText $ws "C38" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-CellText $ws "D38" "12.07"
Set-CellText $ws "E38" "  +1.42%  "

# --- Rows 41/42 swapped position (FraxShare now ranks above TrustWalletToken) ---
Set-CellText $ws "B41" "FraxShare"
Set-CellText $ws "C41" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-CellText $ws "D41" "8.426"
Set-CellText $ws "E41" "  +6.53%  "

Set-CellText $ws "B42" "TrustWalletToken"
Set-CellText $ws "C42" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText $ws "D42" "1.233"
Set-CellText $ws "E42" "  +3.99%  "

# --- Remaining Price (D) / Volume(1h) (E) refreshes ---
Set-CellText $ws "D2" "28.209.21"
Set-CellText $ws "E2" "  +1.99%  "
Set-CellText $ws "D3" "1.799.64"
Set-CellText $ws "E3" "  +3.39%  "
Set-CellText $ws "D4" "1.005"
Set-CellText $ws "E4" "  -0.17%  "
Set-CellText $ws "D5" "338.31"
Set-CellText $ws "E5" "  +2.29%  "
Set-CellText $ws "D6" "1.001"
Set-CellText $ws "D7" "0.4582"
Set-CellText $ws "E7" "  +18.22%  "
Set-CellText $ws "D8" "0.3775"
Set-CellText $ws "E8" "  +13.20%  "
Set-CellText $ws "D9" "45.09"
Set-CellText $ws "E9" "  -0.53%  "
Set-CellText $ws "D10" "0.07629"
Set-CellText $ws "E10" "  +6.55%  "
Set-CellText $ws "D11" "1.146"
Set-CellText $ws "E11" "  +4.33%  "
Set-CellText $ws "D12" "1.002"
Set-CellText $ws "E12" "  -0.14%  "
Set-CellText $ws "D13" "22.39"
Set-CellText $ws "E13" "  +1.34%  "
Set-CellText $ws "D14" "6.332"
Set-CellText $ws "E14" "  +4.06%  "
Set-CellText $ws "D15" "7.503"
Set-CellText $ws "E15" "  +8.10%  "
Set-CellText $ws "D16" "1.801.38"
Set-CellText $ws "E16" "  +3.44%  "
Set-CellText $ws "E17" "  +4.24%  "
Set-CellText $ws "D18" "0.06745"
Set-CellText $ws "E18" "  +2.36%  "
Set-CellText $ws "D19" "81.31"
Set-CellText $ws "E19" "  +4.11%  "
Set-CellText $ws "E20" "  -0.06%  "
Set-CellText $ws "D21" "17.42"
Set-CellText $ws "E21" "  +5.18%  "
Set-CellText $ws "E22" "  +4.24%  "
Set-CellText $ws "D23" "28.238.08"
Set-CellText $ws "E23" "  +2.01%  "
Set-CellText $ws "D24" "11.87"
Set-CellText $ws "E24" "  +3.34%  "
Set-CellText $ws "D25" "2.421"
Set-CellText $ws "E25" "  +1.40%  "
Set-CellText $ws "D26" "20.72"
Set-CellText $ws "E26" "  +5.83%  "
Set-CellText $ws "D27" "152.02"
Set-CellText $ws "E27" "  -2.05%  "
Set-CellText $ws "D28" "2.361"
Set-CellText $ws "E28" "  +5.11%  "
Set-CellText $ws "D29" "2.005.45"
Set-CellText $ws "E29" "  +3.40%  "
Set-CellText $ws "D30" "132.86"
Set-CellText $ws "E30" "  +3.83%  "
Set-CellText $ws "E31" "  -1.42%  "
Set-CellText $ws "D32" "4.030"
Set-CellText $ws "E32" "  +0.33%  "
Set-CellText $ws "D33" "0.09515"
Set-CellText $ws "E33" "  +9.65%  "
Set-CellText $ws "D34" "5.830"
Set-CellText $ws "E34" "  +1.64%  "
Set-CellText $ws "D35" "0.2302"
Set-CellText $ws "E35" "  +10.75%  "
Set-CellText $ws "D36" "0.06347"
Set-CellText $ws "E36" "  +5.63%  "
Set-CellText $ws "D39" "5.246"
Set-CellText $ws "E39" "  +3.71%  "
Set-CellText $ws "D40" "0.6602"
Set-CellText $ws "E40" "  +3.03%  "
Set-CellText $ws "D43" "1.485"
Set-CellText $ws "E43" "  -2.15%  "
Set-CellText $ws "D44" "14.19"
Set-CellText $ws "E44" "  +4.84%  "
Set-CellText $ws "D45" "1.001"
Set-CellText $ws "D46" "3.865"
Set-CellText $ws "E46" "  +1.97%  "
Set-CellText $ws "E47" "  +2.94%  "
Set-CellText $ws "D48" "130.62"
Set-CellText $ws "E48" "  +4.08%  "
Set-CellText $ws "D49" "2.031"
Set-CellText $ws "E49" "  +3.42%  "
Set-CellText $ws "D50" "0.07141"
Set-CellText $ws "D51" "1.165"
Set-CellText $ws "E51" "  +2.01%  "
